{"js": "// The document opens with a centered \"Featured Images\" heading followed\n// by two identical, empty, centered paragraphs (same 48pt Helvetica /\n// underline run formatting) before the body text begins. The edit\n// collapses that to a single empty spacer paragraph by deleting one of\n// the two duplicate empty paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Featured Images\" heading paragraph, then delete the first\n// empty paragraph that immediately follows it (of the pair of empty\n// centered paragraphs under the heading).\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Featured Images\") {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex !== -1 && headingIndex + 2 < paragraphs.items.length) {\n  const first = paragraphs.items[headingIndex + 1];\n  const second = paragraphs.items[headingIndex + 2];\n  if (first.text === \"\" && second.text === \"\") {\n    first.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The document opens with a centered \"Featured Images\" heading followed\n# by two identical, empty, centered paragraphs (same 48pt Helvetica /\n# underline run formatting) before the body text begins. The edit\n# collapses that to a single empty spacer paragraph by deleting one of\n# the two duplicate empty paragraphs.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Featured Images\" heading paragraph (Range.Text carries a\n# trailing paragraph mark, so trim before comparing).\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq \"Featured Images\") {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -ge 1 -and ($headingIndex + 2) -le $d.Paragraphs.Count) {\n    $first = $d.Paragraphs($headingIndex + 1)\n    $second = $d.Paragraphs($headingIndex + 2)\n    if ($first.Range.Text.TrimEnd() -eq \"\" -and $second.Range.Text.TrimEnd() -eq \"\") {\n        $first.Range.Delete()\n    }\n}\n"}
